$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.006.59'
$ws.Range('E2').Value = '  +2.55%  '
$ws.Range('D3').Value = '2.055.33'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '229.90'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('D7').Value = '58.20'
$ws.Range('E7').Value = '  +6.80%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.386'
$ws.Range('E9').Value = '  +2.81%  '
$ws.Range('D10').Value = '0.0808'
$ws.Range('E10').Value = '  +3.15%  '
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').Value = '2.360.53'
$ws.Range('E12').Value = '  +2.01%  '
$ws.Range('D13').Value = '14.59'
$ws.Range('E13').Value = '  +3.32%  '
$ws.Range('D14').Value = '20.68'
$ws.Range('E14').Value = '  +2.40%  '
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').Value = '2.052.84'
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D18').Value = '37.900.20'
$ws.Range('E18').Value = '  +2.42%  '
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').Value = '69.78'
$ws.Range('D21').Value = '0.0₃0831'
$ws.Range('E21').Value = '  +1.79%  '
$ws.Range('D22').Value = '224.93'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('E25').Value = '  +3.47%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '166.36'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '9.28'
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('E28').Value = '  +7.59%  '
$ws.Range('D29').Value = '19.04'
$ws.Range('E29').Value = '  +1.77%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  +1.97%  '
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = '2.03'
$ws.Range('E35').Value = '  +10.12%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +13.07%  '
$ws.Range('E38').Value = '  +5.32%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').Value = '98.12'
$ws.Range('E40').Value = '  +3.31%  '
$ws.Range('D41').Value = '0.0219'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('D42').Value = '1.484.74'
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('E43').Value = '  +3.35%  '
$ws.Range('D44').Value = '0.0936'
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('D45').Value = '16.70'
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = '4.17'
$ws.Range('E46').Value = '  +17.79%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = '1.13'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').Value = '2.96'
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('D50').Value = '7.04'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('D51').Value = '2.248.55'
$ws.Range('E51').Value = '  +2.07%  '
